$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": stamp the latest handoff generation
# timestamps for the be004eef... row (row 7) across the Overview summary
# sheet and the per-locale detail sheets.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-25 12:43:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-25 12:43:45"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-25 12:43:49"
